# Applies the "Adding labs with a correct ic" edit:
#  - Blood pressure reading changes from "123/78" to "124/79" (B1)
#  - Several numeric lab values throughout the sheet are updated
#  - Row heights for the first 4 (header-ish) rows are normalized to 15.75
#  - View/selection state is updated to reflect scrolling near the bottom of the sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Blood pressure text value ---
$ws.Range("B1").Value = "124/79"

# --- Row heights (rows 1-4 become uniform 15.75, matching rest of sheet) ---
$ws.Rows.Item(1).RowHeight = 15.75
$ws.Rows.Item(2).RowHeight = 15.75
$ws.Rows.Item(3).RowHeight = 15.75
$ws.Rows.Item(4).RowHeight = 15.75

# --- Numeric value updates ---
$ws.Range("B4").Value = 98.5

$ws.Range("C7").Value = 144.2
$ws.Range("C8").Value = 32.4
$ws.Range("C9").Value = 107
$ws.Range("C10").Value = 7.38
$ws.Range("C11").Value = 41
$ws.Range("C14").Value = 301
$ws.Range("C15").Value = 45
$ws.Range("C18").Value = 41
$ws.Range("C20").Value = 43
$ws.Range("C21").Value = 48
$ws.Range("C23").Value = 2.3
$ws.Range("C24").Value = 292

$ws.Range("B26").Value = 5346
$ws.Range("B27").Value = 75

$ws.Range("B29").Value = 0.1
$ws.Range("B30").Value = 3.5

$ws.Range("B35").Value = 5413

$ws.Range("B38").Value = 3.4
$ws.Range("B39").Value = 49
$ws.Range("B40").Value = 127

$ws.Range("B42").Value = 38.6
$ws.Range("B43").Value = 12.8

$ws.Range("B46").Value = 6.5
$ws.Range("B47").Value = 4.4
$ws.Range("B48").Value = 516

$ws.Range("B50").Value = 1.211

$ws.Range("B54").Value = 245

$ws.Range("B57").Value = 21

$ws.Range("B59").Value = 20.6

$ws.Range("B61").Value = 77

$ws.Range("B65").Value = 0.66
$ws.Range("B66").Value = 0.125

# --- View state: scroll near the bottom of the sheet and select B66 ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 249
$ws.Range("B66").Select()
